# Atualização de bases das ligas, do dia: 24-02-2024 às 23:13
#
# The source data rows got re-matched to the correct fixtures/odds.
# Concretely, the editable columns B..AC (id, teams, score, odds, ...) of
# certain rows were re-shuffled among themselves while column A (the
# sequential row index) stayed put. We read all the "B:AC" row blocks we
# need to touch FIRST (so the permutation can't clobber data we still need
# to read), and then write them back out according to the new mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return $ws.Range("B" + $row + ":AC" + $row).Value()
}

function Set-RowData($row, $data) {
    $ws.Range("B" + $row + ":AC" + $row).Value = $data
}

# Each inner list is one independent permutation group: row N is overwritten
# with the B:AC content that row M used to hold (before any writes happen).
$groups = @(
    @(425, 426),
    @(470, 473),
    @(497, 499, 500, 498)
)

foreach ($group in $groups) {
    # Snapshot the current ("before") B:AC content of every row in the group.
    $snapshots = @{}
    foreach ($row in $group) {
        $snapshots[$row] = Get-RowData $row
    }

    # group[i] receives the snapshot that used to belong to group[i+1]
    # (wrapping around), i.e. a cyclic rotation of the captured data.
    $count = $group.Count
    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $group[$i]
        $srcRow = $group[($i + 1) % $count]
        Set-RowData $destRow $snapshots[$srcRow]
    }
}
